$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4 (shifts the header/data rows down by one, and copies
# formatting down from row 3, the blank spacer row above the insertion point)
$ws.Rows.Item(4).Insert()

# Set the new footnote text in A3 (the pre-existing blank spacer row, style already 3)
$ws.Range("A3").Value = "Outliers were removed prior to data analysis. One outliers from SY1140A."

# Update the active selection to A3
$ws.Range("A3").Select()
